$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Drop the now-unused trailing rows (old duplicate per-month rows 6-9)
#    and the now-unused trailing column Q, shrinking the used range to A1:P5.
$ws.Range("A9:Q9").EntireRow.Delete()
$ws.Range("A8:Q8").EntireRow.Delete()
$ws.Range("A7:Q7").EntireRow.Delete()
$ws.Range("A6:Q6").EntireRow.Delete()
$ws.Range("Q1:Q9").EntireColumn.Delete()

# 2. Rewrite the header row with the new labels (left to right).
$ws.Range("B1").Value = "负责人"
$ws.Range("C1").Value = "去年实绩"
$ws.Range("D1").Value = "目标值"
$ws.Range("E1").Value = "1月"
$ws.Range("F1").Value = "2月"
$ws.Range("G1").Value = "3月"
$ws.Range("H1").Value = "4月"
$ws.Range("I1").Value = "5月"
$ws.Range("J1").Value = "6月"
$ws.Range("K1").Value = "7月"
$ws.Range("L1").Value = "8月"
$ws.Range("M1").Value = "9月"
$ws.Range("N1").Value = "10月"
$ws.Range("O1").Value = "11月"
$ws.Range("P1").Value = "12月"

# 3. Rewrite the metric-name column (A) for the 4 data rows.
$ws.Range("A2").Value = "1580高牌号硅钢总体一次投料合格率"
$ws.Range("A3").Value = "1580高牌号硅钢板形一次投料合格率"
$ws.Range("A4").Value = "1580高牌号硅钢表面一次投料合格率"
$ws.Range("A5").Value = "1580高牌号硅钢成份性能一次投料合格率"

# 4. Owner column (B) and the numeric data (C: 去年实绩, D: 目标值, E: 1月, F: 2月).
$ws.Range("B2").Value = "王宇阳"
$ws.Range("C2").Value = 62.85
$ws.Range("D2").Value = 68.88
$ws.Range("E2").Value = 45.63
$ws.Range("F2").Value = 58.61

$ws.Range("B3").Value = "王宇阳"
$ws.Range("C3").Value = 92.76000000000001
$ws.Range("D3").Value = 94.77
$ws.Range("E3").Value = 84.78
$ws.Range("F3").Value = 86.41

$ws.Range("B4").Value = "王宇阳"
$ws.Range("C4").Value = 88.20999999999999
$ws.Range("D4").Value = 91.65000000000001
$ws.Range("E4").Value = 84.93000000000001
$ws.Range("F4").Value = 84.06999999999999

$ws.Range("B5").Value = "王宇阳"
$ws.Range("C5").Value = 71.12
$ws.Range("D5").Value = 77.84999999999999
$ws.Range("E5").Value = 58.91
$ws.Range("F5").Value = 70.34
